$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 329-335 with new field values ---
# Row 329
$ws.Cells.Item(329, 4).Value = 44595
$ws.Cells.Item(329, 10).Value = 900
$ws.Cells.Item(329, 11).Value = 8000
$ws.Cells.Item(329, 13).Value = 8822
$ws.Cells.Item(329, 15).Value = "Región Metropolitana"
$ws.Cells.Item(329, 16).Value = 176

# Row 330
$ws.Cells.Item(330, 4).Value = 44595
$ws.Cells.Item(330, 10).Value = 180
$ws.Cells.Item(330, 11).Value = 8000
$ws.Cells.Item(330, 12).Value = 10000
$ws.Cells.Item(330, 13).Value = 9333
$ws.Cells.Item(330, 16).Value = 187

# Row 331
$ws.Cells.Item(331, 4).Value = 44335
$ws.Cells.Item(331, 11).Value = 9000
$ws.Cells.Item(331, 12).Value = 10000
$ws.Cells.Item(331, 13).Value = 9425
$ws.Cells.Item(331, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(331, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(331, 16).Value = 188
$ws.Cells.Item(331, 17).Value = 50

# Row 332
$ws.Cells.Item(332, 4).Value = 44552
$ws.Cells.Item(332, 10).Value = 400
$ws.Cells.Item(332, 11).Value = 7000
$ws.Cells.Item(332, 12).Value = 8000
$ws.Cells.Item(332, 13).Value = 7425
$ws.Cells.Item(332, 16).Value = 148

# Row 333
$ws.Cells.Item(333, 4).Value = 44552
$ws.Cells.Item(333, 10).Value = 400
$ws.Cells.Item(333, 11).Value = 7000
$ws.Cells.Item(333, 12).Value = 8000
$ws.Cells.Item(333, 13).Value = 7425
$ws.Cells.Item(333, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(333, 15).Value = "Región del Maule"
$ws.Cells.Item(333, 16).Value = 124
$ws.Cells.Item(333, 17).Value = 60

# Row 334
$ws.Cells.Item(334, 4).Value = 44544
$ws.Cells.Item(334, 10).Value = 250
$ws.Cells.Item(334, 11).Value = 8000
$ws.Cells.Item(334, 12).Value = 9000
$ws.Cells.Item(334, 13).Value = 8480
$ws.Cells.Item(334, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(334, 16).Value = 170
$ws.Cells.Item(334, 17).Value = 50

# Row 335
$ws.Cells.Item(335, 10).Value = 800
$ws.Cells.Item(335, 11).Value = 5000
$ws.Cells.Item(335, 12).Value = 6000
$ws.Cells.Item(335, 13).Value = 5562
$ws.Cells.Item(335, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(335, 15).Value = "Limache"
$ws.Cells.Item(335, 16).Value = 111
$ws.Cells.Item(335, 17).Value = 50

# --- Append new rows 336 and 337 ---
# Row 336
$ws.Cells.Item(336, 1).Value = 6
$ws.Cells.Item(336, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(336, 3).Value = "Metropolitana"
$ws.Cells.Item(336, 4).Value = 44160
$ws.Cells.Item(336, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(336, 5).Value = 13
$ws.Cells.Item(336, 6).Value = 100112032
$ws.Cells.Item(336, 7).Value = "Zapallo italiano"
$ws.Cells.Item(336, 8).Value = "Sin especificar"
$ws.Cells.Item(336, 9).Value = "Primera"
$ws.Cells.Item(336, 10).Value = 340
$ws.Cells.Item(336, 11).Value = 4000
$ws.Cells.Item(336, 12).Value = 5000
$ws.Cells.Item(336, 13).Value = 4500
$ws.Cells.Item(336, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(336, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(336, 16).Value = 75
$ws.Cells.Item(336, 17).Value = 60
$ws.Cells.Item(336, 18).Value = "Hortaliza"

# Row 337
$ws.Cells.Item(337, 1).Value = 6
$ws.Cells.Item(337, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(337, 3).Value = "Metropolitana"
$ws.Cells.Item(337, 4).Value = 44160
$ws.Cells.Item(337, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(337, 5).Value = 13
$ws.Cells.Item(337, 6).Value = 100112032
$ws.Cells.Item(337, 7).Value = "Zapallo italiano"
$ws.Cells.Item(337, 8).Value = "Sin especificar"
$ws.Cells.Item(337, 9).Value = "Primera"
$ws.Cells.Item(337, 10).Value = 330
$ws.Cells.Item(337, 11).Value = 4000
$ws.Cells.Item(337, 12).Value = 5000
$ws.Cells.Item(337, 13).Value = 4500
$ws.Cells.Item(337, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(337, 15).Value = "Región del Maule"
$ws.Cells.Item(337, 16).Value = 75
$ws.Cells.Item(337, 17).Value = 60
$ws.Cells.Item(337, 18).Value = "Hortaliza"
